$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header G1 from "compra" to "monto"
$ws.Range("G1").Value = "monto"

# Rename header H1 from "venta" to "tipo_instrumento"
$ws.Range("H1").Value = "tipo_instrumento"

# Add new header I1 "comision"
$ws.Range("I1").Value = "comision"

# Copy style from H1 to I1 so the new header matches the bold/bordered look
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2: monto stays (G2 already 10879382), tipo_instrumento = COMPRA, comision = 0 (was H2)
$ws.Range("H2").Value = "COMPRA"
$ws.Range("I2").Value = 0

# Row 3: monto stays (G3 already 1519368), tipo_instrumento = COMPRA, comision = 0 (was H3)
$ws.Range("H3").Value = "COMPRA"
$ws.Range("I3").Value = 0

# Row 4: monto moves from H4 to G4 (was G4=0, H4=314272851), tipo_instrumento = VENTA, comision = 0
$ws.Range("G4").Value = 314272851
$ws.Range("H4").Value = "VENTA"
$ws.Range("I4").Value = 0

# Row 5: monto moves from H5 to G5 (was G5=0, H5=1519368), tipo_instrumento = VENTA, comision = 0
$ws.Range("G5").Value = 1519368
$ws.Range("H5").Value = "VENTA"
$ws.Range("I5").Value = 0
